# Medolla.xlsx - aggiornamento fino a 6/03
# Appends 3 new rows (245-247) of data to the worksheet, extending the
# existing daily series in columns A:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing date cell (A244) onto
# the new date cells so they keep the same date number format/alignment.
$ws.Range("A244").Copy()
$ws.Range("A245:A247").PasteSpecial(-4122)

# Row 245: 2021-05-03
$ws.Cells.Item(245, 1).Value = 44319
$ws.Cells.Item(245, 2).Value = 3
$ws.Cells.Item(245, 3).Value = 11
$ws.Cells.Item(245, 4).Value = 176.253805479891

# Row 246: 2021-05-04
$ws.Cells.Item(246, 1).Value = 44320
$ws.Cells.Item(246, 2).Value = 0
$ws.Cells.Item(246, 3).Value = 9
$ws.Cells.Item(246, 4).Value = 144.2076590290018

# Row 247: 2021-05-05
$ws.Cells.Item(247, 1).Value = 44321
$ws.Cells.Item(247, 2).Value = 0
$ws.Cells.Item(247, 3).Value = 9
$ws.Cells.Item(247, 4).Value = 144.2076590290018
